# "rename everything in english" - translate remaining German labels to English
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "Balance"          # was "Kontostand"
$ws.Range("A3").Value = "Income"           # was "Einnahmen"
$ws.Range("A4").Value = "Expenses"         # was "Ausgaben"
$ws.Range("A5").Value = "Balance"          # was "Bilanz"
$ws.Range("A7").Value = "CATEGORIES"       # was "Category"
$ws.Range("B9").Value = "SUBCATEGORY"      # was "UNTERKATEGORIE"
$ws.Range("A8").Value = "TOP-CATEGORY"     # was "OBERKATEGORIE"
$ws.Range("A10").Value = "OTHERS"          # was "SONSTIGES"

# Update the selection to match the recorded cursor position
$ws.Range("A10:B10").Select()
